# Update patient record fields in the admission/discharge sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apellidos / Nombres (row 6)
$ws.Range("A6").Value = "MEJÌA"
$ws.Range("C6").Value = "AMPEREZ"
$ws.Range("E6").Value = "MARÌA"
$ws.Range("G6").Value = "ADELA "
$ws.Range("I6").Value = "/201773414"

# Direccion actual (row 8)
$ws.Range("A8").Value = "LOTE 25 MANZ 10 "
$ws.Range("D8").Value = "JUANA DE ARCO Z. 18"

# Fecha de nacimiento / Edad / Lugar de nacimiento (row 12)
$ws.Range("A12").Value = "1996-06-10"
$ws.Range("F12").Value = "22"
$ws.Range("H12").Value = "SAN MIGUEL IXTAHUACÀN/SAN MARCOS"

# Estado civil / Nacionalidad / No. de cedula (row 14)
$ws.Range("A14").Value = "Soltero"
$ws.Range("F14").Value = "GUATEMALTECA"
$ws.Range("H14").Value = "33260606361205"

# Direccion si difiere a la indicada (row 16) - cleared
$ws.Range("A16").Value = ""

# Nombre del padre / madre (row 18)
$ws.Range("A18").Value = "MARIO MEJÌA"
$ws.Range("F18").Value = "CATALINA AMPEREZ"

# En caso de emergencia (row 20)
$ws.Range("A20").Value = ""
$ws.Range("F20").Value = "MADRE"
$ws.Range("J20").Value = "5958-6060"

# Fecha / hora de ingreso (row 24)
$ws.Range("A24").Value = "20/11/2017"
$ws.Range("C24").Value = "12:36:18"
